# Refresh the cryptocurrency price ("D") and 1h volume-change ("E") columns
# on Sheet1 to match the latest scrape, cell by cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain numbers (e.g. "15.60", "0.9987").
# Prefixing with an apostrophe forces Excel to keep them as literal text
# (matching the source data), instead of parsing them into numbers and
# silently dropping significant trailing zeros.

$ws.Range('D2').Value = '26.579.52'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').Value = '1.731.43'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D4').Value = '''0.9987'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''244.99'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D7').Value = '''0.4812'
$ws.Range('E7').Value = '  +1.48%  '
$ws.Range('D8').Value = '''0.2678'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').Value = '''0.06194'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').Value = '1.727.98'
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('D11').Value = '''0.07189'
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('D12').Value = '''15.60'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '''0.6115'
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').Value = '''77.31'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '''0.9991'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '26.570.44'
$ws.Range('D18').Value = '''0.9994'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '''0.000006967'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('D20').Value = '''11.56'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '1.952.43'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('D22').Value = '''4.532'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').Value = '''8.821'
$ws.Range('E23').Value = '  +0.70%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').Value = '''15.37'
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('D27').Value = '''1.785'
$ws.Range('E27').Value = '  +0.73%  '
$ws.Range('D28').Value = '''1.406'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '''107.36'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('D30').Value = '''3.987'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('D31').Value = '''0.08036'
$ws.Range('E31').Value = '  +2.90%  '
$ws.Range('D32').Value = '''3.701'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').Value = '''0.04531'
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('E35').Value = '  +2.87%  '
$ws.Range('D36').Value = '''0.6280'
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').Value = '''2.088'
$ws.Range('E37').Value = '  +7.94%  '
$ws.Range('D38').Value = '''0.9074'
$ws.Range('E38').Value = '  -4.12%  '
$ws.Range('D39').Value = '''2.406'
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').Value = '''1.003'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').Value = '''0.01505'
$ws.Range('E41').Value = '  +1.23%  '
$ws.Range('D42').Value = '''102.50'
$ws.Range('E42').Value = '  -10.40%  '
$ws.Range('D43').Value = '''5.536'
$ws.Range('E43').Value = '  -2.83%  '
$ws.Range('D44').Value = '''0.3894'
$ws.Range('E44').Value = '  +1.48%  '
$ws.Range('D45').Value = '''7.022'
$ws.Range('E45').Value = '  +10.07%  '
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('D47').Value = '''0.05380'
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('D48').Value = '''30.74'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').Value = '''7.834'
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '''1.252'
$ws.Range('E50').Value = '  +2.68%  '
$ws.Range('D51').Value = '''0.3418'
$ws.Range('E51').Value = '  +0.79%  '
